# "Half working 'Note' page" — add a new (blank) slide #3 that holds a
# single, barely-sized text box pre-wired for a Segoe MDL2 Assets glyph
# (the icon itself was never actually typed in, hence "half working").

$p = $ppt.ActivePresentation

# New slide goes after the two existing ones, using the "blank" layout
# (ppLayoutBlank = 12) -- mirrors slideLayouts/slideLayout7.xml.
$s = $p.Slides.Add($p.Slides.Count + 1, 12)

# The shape-id counter is shared per-slide and starts at 2 (1 is the root
# group). The real file's textbox carries id="4", so burn two throw-away
# ids first and remove those placeholders again, leaving a clean,
# single-shape slide.
$junk1 = $s.Shapes.AddTextbox(1, 0, 0, 1, 1)
$junk2 = $s.Shapes.AddTextbox(1, 0, 0, 1, 1)
$junk1.Delete()
$junk2.Delete()

# Shape geometry is in points; the target EMU box is
# off=(238539,357809) ext=(45719,369332) -> divide by 12700 (EMU/pt).
$left   = 238539 / 12700.0
$top    = 357809 / 12700.0
$width  = 45719 / 12700.0
$height = 369332 / 12700.0

$shp = $s.Shapes.AddTextbox(1, $left, $top, $width, $height)
$shp.Name = "กล่องข้อความ 3"

# wrap="square" + spAutoFit
$shp.TextFrame.WordWrap = 1
$shp.TextFrame.AutoSize = 1

# noFill
$shp.Fill.Visible = 0

# Text run pre-set to the Segoe MDL2 Assets icon font; the glyph itself
# was never filled in (hence "half working"), so the run text is empty.
$tr = $shp.TextFrame.TextRange
$tr.Text = ""
$tr.Font.Name = "Segoe MDL2 Assets"
